$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 131, shifting rows 131:167 down to 132:168,
# and copying formatting (e.g. the date style) down from the row above.
$ws.Rows("131").Insert()

# Populate the newly-inserted row 131 with the new weekly record.
$ws.Range("A131").Value = 9
$ws.Range("B131").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C131").Value = "Metropolitana"
$ws.Range("D131").Value = 45093
$ws.Range("E131").Value = 13
$ws.Range("F131").Value = 100112022
$ws.Range("G131").Value = "Arveja Verde"
$ws.Range("H131").Value = "Perfection"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 52
$ws.Range("K131").Value = 40000
$ws.Range("L131").Value = 42000
$ws.Range("M131").Value = 41000
$ws.Range("N131").Value = "$/saco 25 kilos"
$ws.Range("O131").Value = "Provincia de Huasco"
$ws.Range("P131").Value = 1640
$ws.Range("Q131").Value = 25
$ws.Range("R131").Value = "Hortaliza"
